# Template surat selesai: append a "/${tahun_surat_selesai}" placeholder
# right after the existing "${nomor_surat_selesai}" placeholder in the
# "Nomor : 400.14.5.4/${nomor_surat_selesai}" line, so the final text
# reads "Nomor : 400.14.5.4/${nomor_surat_selesai}/${tahun_surat_selesai}".

$d = $word.ActiveDocument

# Locate the unique "${nomor_surat_selesai}" placeholder text.
$r = $d.Content
$r.Find.Execute("`${nomor_surat_selesai`}", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)

# Pull the range's end back one character so it sits right before the
# closing "}" of "${nomor_surat_selesai}" (i.e. right after "nomor_surat_selesai"),
# leaving that final "}" run untouched.
$null = $r.MoveEnd(1, -1)
$r.Collapse(0)

# Insert the new "}/${tahun_surat_selesai" text before the original,
# unchanged closing "}" run — yielding
# "${nomor_surat_selesai}/${tahun_surat_selesai}" overall.
$r.InsertAfter("`}/`${tahun_surat_selesai")

Write-Output "Appended /`${tahun_surat_selesai} after the nomor_surat_selesai placeholder."
